$wb = $excel.ActiveWorkbook

# The most recent existing sheet (index 1) holds the latest snapshot
# ("2025-06-05"). Duplicate it to create the next day's snapshot
# ("2025-06-08"), inserted before it, exactly like Excel's own
# "Move or Copy... (Create a copy)" in front of the source sheet.
$sourceSheet = $wb.Worksheets.Item(1)
$sourceSheet.Copy($sourceSheet)

# The newly inserted copy is now the first sheet; rename it.
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "2025-06-08"

# Restore the previously-active sheet's plain selection (A1:F1), then
# make the new sheet active with the cursor at G1.
$previousSheet = $wb.Worksheets.Item(2)
[void]$previousSheet.Range("A1:F1").Select()

$newSheet.Activate()
[void]$newSheet.Range("G1").Select()
